# Insert a new weekly price record for "Feria Lagunitas de Puerto Montt" (Cilantro)
# as row 456, pushing every subsequent row down by one (old row 456 becomes 457,
# ..., old row 536 becomes the new last row 537). The worksheet's used
# range grows from A1:R536 to A1:R537 accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 456..536 down by one row, creating an empty row 456.
$ws.Rows("456:456").Insert()

# Populate the newly inserted row with the new data point.
$ws.Range("A456").Value = 4
$ws.Range("B456").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C456").Value = "Los Lagos"
$ws.Range("D456").Value = 45209
$ws.Range("E456").Value = 10
$ws.Range("F456").Value = 100112040
$ws.Range("G456").Value = "Cilantro"
$ws.Range("H456").Value = "Sin especificar"
$ws.Range("I456").Value = "Primera"
$ws.Range("J456").Value = 180
$ws.Range("K456").Value = 13000
$ws.Range("L456").Value = 13000
$ws.Range("M456").Value = 13000
$ws.Range("N456").Value = "$/caja 36 atados"
$ws.Range("O456").Value = "Región Metropolitana"
$ws.Range("P456").Value = 361
$ws.Range("Q456").Value = 36
$ws.Range("R456").Value = "Hortaliza"
